$d = $word.ActiveDocument

# 1. Correct "Vidharmi" -> "Vidharma" within the italic sentence about the
#    new sect's adherents.
$r = $d.Content
$found = $r.Find.Execute("Vidharmi", $true, $false, $false, $false, $false, `
                          $true, 1, $false, "Vidharma", 2)

# 2. Word keeps exactly one "_GoBack" bookmark, always relocated to the site
#    of the most recent edit. Move it from wherever it currently sits to
#    right after the word we just corrected.
if ($found) {
    $editEnd = $r.End

    if ($d.Bookmarks.Exists("_GoBack")) {
        $d.Bookmarks("_GoBack").Delete()
    }

    $bmRange = $d.Range($editEnd, $editEnd)
    $d.Bookmarks.Add("_GoBack", $bmRange)
}
